$d = $word.ActiveDocument

# Build a single-paragraph WordprocessingML "InsertXML" payload containing one
# run (optionally with the given <w:rPr>...</w:rPr> fragment) holding $text.
# Using InsertXML (rather than Range.Text / Find&Replace) on the *exact*
# found range is important here: this engine normalizes (drops) sibling
# empty <w:r/> runs in a paragraph whenever Range.Text / Find-replace / Delete
# touch that paragraph, but InsertXML surgically replaces only the targeted
# range and leaves neighboring runs (like a leading empty <w:r/>) untouched.
function New-RunPackageXml([string]$text, [string]$rPrXml) {
    $runXml = "<w:r>" + $rPrXml + "<w:t>" + $text + "</w:t></w:r>"
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
           $runXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Replace every exact occurrence of $old (whole string match) with a single
# run containing $new, preserving the run's formatting ($rPrXml) and any
# sibling empty runs in that paragraph.
function Replace-RunText([string]$old, [string]$new, [string]$rPrXml = "") {
    $searchFrom = 0
    while ($true) {
        $scan = $d.Range($searchFrom, $d.Content.End)
        $found = $scan.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $found) { break }
        $s = $scan.Start
        $e = $scan.End
        $target = $d.Range($s, $e)
        $target.InsertXML((New-RunPackageXml $new $rPrXml))
        $searchFrom = $s + $new.Length
    }
}

# Replace only the Nth (1-based) occurrence of $old, keeping its own
# formatting via $rPrXml.
function Replace-RunTextOccurrence([string]$old, [string]$new, [int]$occurrence, [string]$rPrXml = "") {
    $searchFrom = 0
    $n = 0
    while ($true) {
        $scan = $d.Range($searchFrom, $d.Content.End)
        $found = $scan.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $found) { break }
        $n = $n + 1
        $s = $scan.Start
        $e = $scan.End
        if ($n -eq $occurrence) {
            $target = $d.Range($s, $e)
            $target.InsertXML((New-RunPackageXml $new $rPrXml))
            break
        }
        $searchFrom = $e
    }
}

# Title: 1st occurrence is the plain H1 heading, 2nd is a bold run further down.
Replace-RunTextOccurrence "Play Mother of Horus Free - Exciting Ancient Egyptian Slot" `
                           "Play Mother of Horus Free - Exciting Gameplay and Stunning Design" 1
Replace-RunTextOccurrence "Play Mother of Horus Free - Exciting Ancient Egyptian Slot" `
                           "Play Mother of Horus Free - Exciting Gameplay and Stunning Design" 1 `
                           "<w:rPr><w:b/></w:rPr>"

# "What we like" bullets
Replace-RunText "Enchanting ancient Egyptian theme" "Enchanting design and stunning graphics"
Replace-RunText "Stunningly detailed graphics and visuals" "Exciting bonus features and mini-games"
Replace-RunText "Several bonus features for exciting gameplay" "High potential payouts"
Replace-RunText "Potential for big wins" "Thrilling slot experience"

# "What we don't like" bullets
Replace-RunText "Moderately low RTP at around 95.30%" "Moderately low RTP"
Replace-RunText "Moderate volatility level" "Moderate volatility"

# Closing italic summary paragraph
$oldSummary = "Read our review of Mother of Horus, an exciting online video slot with " + `
              "free play. Enjoy the enchanting ancient Egyptian theme and potential for big wins."
$newSummary = "Experience the thrill of Mother of Horus with its enchanting design and " + `
              "exciting gameplay. Play for free now."
Replace-RunText $oldSummary $newSummary "<w:rPr><w:i/></w:rPr>"
